$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031523182197974
$ws.Range("D2").Value = 1.035472440590324
$ws.Range("E2").Value = 1.041517427706352
$ws.Range("F2").Value = 1.053563000024011
$ws.Range("I2").Value = 1.038601456581881
$ws.Range("J2").Value = 1.036658439809023
$ws.Range("K2").Value = 1.038268947063911
$ws.Range("L2").Value = 1.044296718287106
$ws.Range("M2").Value = 1.056308623882265
$ws.Range("N2").Value = 1.016216545043877
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032295380322781
$ws.Range("D3").Value = 1.0360385339644
$ws.Range("E3").Value = 1.042240567107563
$ws.Range("F3").Value = 1.054452511640044
$ws.Range("I3").Value = 1.038787632638374
$ws.Range("J3").Value = 1.037073451301886
$ws.Range("K3").Value = 1.03864508943831
$ws.Range("L3").Value = 1.044830731313526
$ws.Range("M3").Value = 1.057011009207023
$ws.Range("N3").Value = 1.016354433564979
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.032795615208217
$ws.Range("D4").Value = 1.036405233761601
$ws.Range("E4").Value = 1.042709422357419
$ws.Range("F4").Value = 1.055029293483402
$ws.Range("I4").Value = 1.038907047252096
$ws.Range("J4").Value = 1.03734187632666
$ws.Range("K4").Value = 1.03888817409664
$ws.Range("L4").Value = 1.045176537192372
$ws.Range("M4").Value = 1.05746606878048
$ws.Range("N4").Value = 1.016443599766802
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033006048847242
$ws.Range("D5").Value = 1.036559488173092
$ws.Range("E5").Value = 1.042906751182873
$ws.Range("F5").Value = 1.05527205978837
$ws.Range("I5").Value = 1.03895699603791
$ws.Range("J5").Value = 1.037454693362121
$ws.Range("K5").Value = 1.038990292682743
$ws.Range("L5").Value = 1.045321975621784
$ws.Range("M5").Value = 1.057657510328673
$ws.Range("N5").Value = 1.01648107116505
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033041389438484
$ws.Range("D6").Value = 1.036585393615662
$ws.Range("E6").Value = 1.042939896527572
$ws.Range("F6").Value = 1.055312838101834
$ws.Range("I6").Value = 1.038965367791209
$ws.Range("J6").Value = 1.03747363411023
$ws.Range("K6").Value = 1.03900743445408
$ws.Range("L6").Value = 1.045346398944821
$ws.Range("M6").Value = 1.057689662042358
$ws.Range("N6").Value = 1.016487361939269
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032798426505043
$ws.Range("D7").Value = 1.036407294550269
$ws.Range("E7").Value = 1.042712058204655
$ws.Range("F7").Value = 1.05503253621267
$ws.Range("I7").Value = 1.038907715665857
$ws.Range("J7").Value = 1.037343383908351
$ws.Range("K7").Value = 1.038889538903365
$ws.Range("L7").Value = 1.045178480307477
$ws.Range("M7").Value = 1.057468626305555
$ws.Range("N7").Value = 1.016444100517326
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031784030832346
$ws.Range("D8").Value = 1.035663670670884
$ws.Range("E8").Value = 1.041761621149111
$ws.Range("F8").Value = 1.053863363202445
$ws.Range("I8").Value = 1.03866459332563
$ws.Range("J8").Value = 1.036798717964893
$ws.Range("K8").Value = 1.038396128442641
$ws.Range("L8").Value = 1.044477134976894
$ws.Range("M8").Value = 1.056545879338463
$ws.Range("N8").Value = 1.016263156622754
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.030000990274164
$ws.Range("D9").Value = 1.0343564512138
$ws.Range("E9").Value = 1.040094081340704
$ws.Range("F9").Value = 1.051812476763757
$ws.Range("I9").Value = 1.038228147867104
$ws.Range("J9").Value = 1.035838121340164
$ws.Range("K9").Value = 1.037524399796866
$ws.Range("L9").Value = 1.043243361284975
$ws.Range("M9").Value = 1.054924320573334
$ws.Range("N9").Value = 1.015943894657711
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028815400006235
$ws.Range("D10").Value = 1.033487192624611
$ws.Range("E10").Value = 1.03898737911528
$ws.Range("F10").Value = 1.050451625977938
$ws.Range("I10").Value = 1.037931834137077
$ws.Range("J10").Value = 1.035197244562402
$ws.Range("K10").Value = 1.036941793253258
$ws.Range("L10").Value = 1.042422337981934
$ws.Range("M10").Value = 1.053846372954437
$ws.Range("N10").Value = 1.015730800878566
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028302785363743
$ws.Range("D11").Value = 1.03311134401623
$ws.Range("E11").Value = 1.038509373047754
$ws.Range("F11").Value = 1.049863908247345
$ws.Range("I11").Value = 1.037802268751735
$ws.Range("J11").Value = 1.034919640301177
$ws.Range("K11").Value = 1.03668918943804
$ws.Range("L11").Value = 1.042067198800198
$ws.Range("M11").Value = 1.053380365597902
$ws.Range("N11").Value = 1.01563847474976
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028112492557862
$ws.Range("D12").Value = 1.032971821035224
$ws.Range("E12").Value = 1.038332003080292
$ws.Range("F12").Value = 1.049645837144162
$ws.Range("I12").Value = 1.037753954008915
$ws.Range("J12").Value = 1.034816511932354
$ws.Range("K12").Value = 1.036595312654923
$ws.Range("L12").Value = 1.041935341318129
$ws.Range("M12").Value = 1.053207384386017
$ws.Range("N12").Value = 1.015604172894942
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028153305772589
$ws.Range("D13").Value = 1.033001745360237
$ws.Range("E13").Value = 1.038370041237835
$ws.Range("F13").Value = 1.049692603538865
$ws.Range("I13").Value = 1.037764326196659
$ws.Range("J13").Value = 1.034838633920582
$ws.Range("K13").Value = 1.036615451716479
$ws.Range("L13").Value = 1.041963622590598
$ws.Range("M13").Value = 1.053244484224776
$ws.Range("N13").Value = 1.015611531105709
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.0282870533424
$ws.Range("D14").Value = 1.033099809280147
$ws.Range("E14").Value = 1.038494707850992
$ws.Range("F14").Value = 1.049845877647784
$ws.Range("I14").Value = 1.037798278878205
$ws.Range("J14").Value = 1.034911115951773
$ws.Range("K14").Value = 1.036681430535302
$ws.Range("L14").Value = 1.042056298247477
$ws.Range("M14").Value = 1.053366064567171
$ws.Range("N14").Value = 1.015635639504598
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028369474941195
$ws.Range("D15").Value = 1.033160240875052
$ws.Range("E15").Value = 1.038571543343972
$ws.Range("F15").Value = 1.049940345880041
$ws.Range("I15").Value = 1.037819173307809
$ws.Range("J15").Value = 1.034955772735394
$ws.Range("K15").Value = 1.036722075878105
$ws.Range("L15").Value = 1.042113406380322
$ws.Range("M15").Value = 1.053440989478273
$ws.Range("N15").Value = 1.015650492462411
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028849436590827
$ws.Range("D16").Value = 1.033512148125295
$ws.Range("E16").Value = 1.039019128325511
$ws.Range("F16").Value = 1.050490663487223
$ws.Range("I16").Value = 1.037940406499394
$ws.Range("J16").Value = 1.035215666238569
$ws.Range("K16").Value = 1.036958550866707
$ws.Range("L16").Value = 1.042445915312176
$ws.Range("M16").Value = 1.053877316308328
$ws.Range("N16").Value = 1.01573692713778
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02915070696475
$ws.Range("D17").Value = 1.033733037781615
$ws.Range("E17").Value = 1.039300210120236
$ws.Range("F17").Value = 1.0508362768824
$ws.Range("I17").Value = 1.038016116249614
$ws.Range("J17").Value = 1.03537866460556
$ws.Range("K17").Value = 1.037106797589706
$ws.Range("L17").Value = 1.042654589343767
$ws.Range("M17").Value = 1.054151214976426
$ws.Range("N17").Value = 1.015791130895387
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029326505480147
$ws.Range("D18").Value = 1.033861931545457
$ws.Range("E18").Value = 1.039464276343355
$ws.Range("F18").Value = 1.051038015739451
$ws.Range("I18").Value = 1.038060154877373
$ws.Range("J18").Value = 1.035473729022516
$ws.Range("K18").Value = 1.037193235338701
$ws.Range("L18").Value = 1.04277634090763
$ws.Range("M18").Value = 1.054311047780813
$ws.Range("N18").Value = 1.015822741681355
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029386460489338
$ws.Range("D19").Value = 1.033905889846543
$ws.Range("E19").Value = 1.039520238295248
$ws.Range("F19").Value = 1.051106828584992
$ws.Range("I19").Value = 1.038075150262034
$ws.Range("J19").Value = 1.03550614183493
$ws.Range("K19").Value = 1.037222702900446
$ws.Range("L19").Value = 1.042817861024634
$ws.Range("M19").Value = 1.054365558852043
$ws.Range("N19").Value = 1.015833519221199
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029118375981761
$ws.Range("D20").Value = 1.033709332961406
$ws.Range("E20").Value = 1.03927004069437
$ws.Range("F20").Value = 1.050799180469844
$ws.Range("I20").Value = 1.038008005880291
$ws.Range("J20").Value = 1.035361177424881
$ws.Range("K20").Value = 1.037090895425682
$ws.Range("L20").Value = 1.042632196908762
$ws.Range("M20").Value = 1.05412182074529
$ws.Range("N20").Value = 1.015785315895101
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028247664824559
$ws.Range("D21").Value = 1.033070929568029
$ws.Range("E21").Value = 1.038457991579378
$ws.Range("F21").Value = 1.049800735782151
$ws.Range("I21").Value = 1.037788285853398
$ws.Range("J21").Value = 1.034889772170307
$ws.Range("K21").Value = 1.036662002745455
$ws.Range("L21").Value = 1.042029005989977
$ws.Range("M21").Value = 1.053330259006526
$ws.Range("N21").Value = 1.015628540394176
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027700880991938
$ws.Range("D22").Value = 1.032670026519483
$ws.Range("E22").Value = 1.037948482499082
$ws.Range("F22").Value = 1.049174325581329
$ws.Range("I22").Value = 1.037649049607318
$ws.Range("J22").Value = 1.034593301856929
$ws.Range("K22").Value = 1.036392061053289
$ws.Range("L22").Value = 1.041650086593313
$ws.Range("M22").Value = 1.052833237070878
$ws.Range("N22").Value = 1.015529924399315
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027990677620708
$ws.Range("D23").Value = 1.032882506068725
$ws.Range("E23").Value = 1.03821848191305
$ws.Range("F23").Value = 1.049506268657707
$ws.Range("I23").Value = 1.037722964368406
$ws.Range("J23").Value = 1.034750473446503
$ws.Range("K23").Value = 1.036535188343518
$ws.Range("L23").Value = 1.04185092707828
$ws.Range("M23").Value = 1.053096654219923
$ws.Range("N23").Value = 1.01558220671725
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029132984737207
$ws.Range("D24").Value = 1.033720043987132
$ws.Range("E24").Value = 1.039283672601232
$ws.Range("F24").Value = 1.050815942283589
$ws.Range("I24").Value = 1.038011670982877
$ws.Range("J24").Value = 1.035369079159467
$ws.Range("K24").Value = 1.037098081028322
$ws.Range("L24").Value = 1.042642314976681
$ws.Range("M24").Value = 1.054135102509779
$ws.Range("N24").Value = 1.015787943460013
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03046140995661
$ws.Range("D25").Value = 1.034694015292553
$ws.Range("E25").Value = 1.040524309094639
$ws.Range("F25").Value = 1.052341559844892
$ws.Range("I25").Value = 1.038341926187996
$ws.Range("J25").Value = 1.036086548593023
$ws.Range("K25").Value = 1.037750024433113
$ws.Range("L25").Value = 1.043562064861953
$ws.Range("M25").Value = 1.055342995297056
$ws.Range("N25").Value = 1.016026478004553
